# Update recomputed TPM-based statistics in the LR-pairs data sheet
# (App -> Ncstn ligand-receptor pair table), per new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "G2" = 89.38217433333334
    "H2" = 268.146523
    "I2" = 0.2143552015363441
    "J2" = 0.2175965347165783
    "M2" = 17.723347
    "N2" = 53.170041
    "O2" = 0.1083389314942055
    "P2" = 0.1121884745845309
    "Q2" = 1584.15129132416
    "R2" = 14257.36162191744
    "S2" = 0.02322301349467259
    "T2" = 0.02441182330473285
    "G3" = 89.38217433333334
    "H3" = 268.146523
    "I3" = 0.2143552015363441
    "J3" = 0.2175965347165783
    "O3" = 0.193467093096278
    "P3" = 0.2003414447366573
    "Q3" = 2828.910541485197
    "R3" = 25460.19487336677
    "S3" = 0.04147067773130332
    "T3" = 0.04359360413480951
    "G4" = 89.38217433333334
    "H4" = 268.146523
    "I4" = 0.2143552015363441
    "J4" = 0.2175965347165783
    "M4" = 47.45519633333333
    "N4" = 142.365589
    "O4" = 0.2900832029413559
    "P4" = 0.3003905575931054
    "Q4" = 4241.648631688561
    "R4" = 38174.83768519705
    "S4" = 0.06218084342880254
    "T4" = 0.06536394439384048
    "G5" = 89.38217433333334
    "H5" = 268.146523
    "I5" = 0.2143552015363441
    "J5" = 0.2175965347165783
    "M5" = 16.8400505
    "N5" = 33.680101
    "O5" = 0.102939533795646
    "P5" = 0.07106481552351887
    "Q5" = 1505.200329573137
    "R5" = 9031.201977438823
    "S5" = 0.022065624512823
    "T5" = 0.01546345759819061
    "G6" = 89.38217433333334
    "H6" = 268.146523
    "I6" = 0.2143552015363441
    "J6" = 0.2175965347165783
    "M6" = 49.92347333333333
    "N6" = 149.77042
    "O6" = 0.3051712386725145
    "P6" = 0.3160147075621876
    "Q6" = 4462.268596805518
    "R6" = 40160.41737124966
    "S6" = 0.06541504236874261
    "T6" = 0.0687637052850049
    "I7" = 0.2934277926151677
    "J7" = 0.2978648075949286
    "M7" = 17.723347
    "N7" = 53.170041
    "O7" = 0.1083389314942055
    "P7" = 0.1121884745845309
    "Q7" = 2168.522215696747
    "R7" = 19516.69994127073
    "S7" = 0.03178965352263058
    "T7" = 0.03341699839648984
    "I8" = 0.2934277926151677
    "J8" = 0.2978648075949286
    "O8" = 0.193467093096278
    "P8" = 0.2003414447366573
    "S8" = 0.056768622070914
    "T8" = 0.05967466588977444
    "I9" = 0.2934277926151677
    "J9" = 0.2978648075949286
    "M9" = 47.45519633333333
    "N9" = 142.365589
    "O9" = 0.2900832029413559
    "P9" = 0.3003905575931054
    "Q9" = 5806.332601798304
    "R9" = 52256.99341618473
    "S9" = 0.08511847391381976
    "T9" = 0.08947577564080364
    "I10" = 0.2934277926151677
    "J10" = 0.2978648075949286
    "M10" = 16.8400505
    "N10" = 33.680101
    "O10" = 0.102939533795646
    "P10" = 0.07106481552351887
    "Q10" = 2060.447364863145
    "R10" = 12362.68418917887
    "S10" = 0.03020532017449085
    "T10" = 0.02116770760268204
    "I11" = 0.2934277926151677
    "J11" = 0.2978648075949286
    "M11" = 49.92347333333333
    "N11" = 149.77042
    "O11" = 0.3051712386725145
    "P11" = 0.3160147075621876
    "Q11" = 6108.336140350775
    "R11" = 54975.02526315697
    "S11" = 0.08954572293331242
    "T11" = 0.09412966006517862
    "G12" = 90.33462533333334
    "H12" = 271.003876
    "I12" = 0.2166393574945233
    "J12" = 0.2199152301234996
    "M12" = 17.723347
    "N12" = 53.170041
    "O12" = 0.1083389314942055
    "P12" = 0.1121884745845309
    "Q12" = 1601.031910897657
    "R12" = 14409.28719807892
    "S12" = 0.02347047651054785
    "T12" = 0.0246719542054615
    "G13" = 90.33462533333334
    "H13" = 271.003876
    "I13" = 0.2166393574945233
    "J13" = 0.2199152301234996
    "O13" = 0.193467093096278
    "P13" = 0.2003414447366573
    "Q13" = 2859.055239734535
    "R13" = 25731.49715761081
    "S13" = 0.04191258674471079
    "T13" = 0.04405813492253637
    "G14" = 90.33462533333334
    "H14" = 271.003876
    "I14" = 0.2166393574945233
    "J14" = 0.2199152301234996
    "M14" = 47.45519633333333
    "N14" = 142.365589
    "O14" = 0.2900832029413559
    "P14" = 0.3003905575931054
    "Q14" = 4286.84738089144
    "R14" = 38581.62642802297
    "S14" = 0.06284343870516873
    "T14" = 0.06606045860001414
    "G15" = 90.33462533333334
    "H15" = 271.003876
    "I15" = 0.2166393574945233
    "J15" = 0.2199152301234996
    "M15" = 16.8400505
    "N15" = 33.680101
    "O15" = 0.102939533795646
    "P15" = 0.07106481552351887
    "Q15" = 1521.239652511913
    "R15" = 9127.437915071476
    "S15" = 0.02230075446227451
    "T15" = 0.0156282352595387
    "G16" = 90.33462533333334
    "H16" = 271.003876
    "I16" = 0.2166393574945233
    "J16" = 0.2199152301234996
    "M16" = 49.92347333333333
    "N16" = 149.77042
    "O16" = 0.3051712386725145
    "P16" = 0.3160147075621876
    "Q16" = 4509.818258905325
    "R16" = 40588.36433014792
    "S16" = 0.06611210107182135
    "T16" = 0.06949644713594892
    "G17" = 18.634161
    "H17" = 37.268322
    "I17" = 0.0446882095496985
    "J17" = 0.03024263611988591
    "M17" = 17.723347
    "N17" = 53.170041
    "O17" = 0.1083389314942055
    "P17" = 0.1121884745845309
    "Q17" = 330.259701456867
    "R17" = 1981.558208741202
    "S17" = 0.004841472873003485
    "T17" = 0.003392875213705038
    "G18" = 18.634161
    "H18" = 37.268322
    "I18" = 0.0446882095496985
    "J18" = 0.03024263611988591
    "O18" = 0.193467093096278
    "P18" = 0.2003414447366573
    "Q18" = 589.7638413677919
    "R18" = 3538.583048206752
    "S18" = 0.0086456979972575
    "T18" = 0.00605885341290296
    "G19" = 18.634161
    "H19" = 37.268322
    "I19" = 0.0446882095496985
    "J19" = 0.03024263611988591
    "M19" = 47.45519633333333
    "N19" = 142.365589
    "O19" = 0.2900832029413559
    "P19" = 0.3003905575931054
    "Q19" = 884.287768761943
    "R19" = 5305.726612571658
    "S19" = 0.01296329895989103
    "T19" = 0.00908460232713792
    "G20" = 18.634161
    "H20" = 37.268322
    "I20" = 0.0446882095496985
    "J20" = 0.03024263611988591
    "M20" = 16.8400505
    "N20" = 33.680101
    "O20" = 0.102939533795646
    "P20" = 0.07106481552351887
    "Q20" = 313.8002122651305
    "R20" = 1255.200849060522
    "S20" = 0.004600183457208099
    "T20" = 0.002149187356804601
    "G21" = 18.634161
    "H21" = 37.268322
    "I21" = 0.0446882095496985
    "J21" = 0.03024263611988591
    "M21" = 49.92347333333333
    "N21" = 149.77042
    "O21" = 0.3051712386725145
    "P21" = 0.3160147075621876
    "Q21" = 930.28203977254
    "R21" = 5581.69223863524
    "S21" = 0.01363755626233838
    "T21" = 0.009557117809335398
    "G22" = 96.27664699999998
    "H22" = 288.829941
    "I22" = 0.2308894388042666
    "J22" = 0.2343807914451077
    "M22" = 17.723347
    "N22" = 53.170041
    "O22" = 0.1083389314942055
    "P22" = 0.1121884745845309
    "Q22" = 1706.344422777509
    "R22" = 15357.09980499758
    "S22" = 0.02501431509335099
    "T22" = 0.0262948234641417
    "G23" = 96.27664699999998
    "H23" = 288.829941
    "I23" = 0.2308894388042666
    "J23" = 0.2343807914451077
    "O23" = 0.193467093096278
    "P23" = 0.2003414447366573
    "Q23" = 3047.11788036665
    "R23" = 27424.06092329985
    "S23" = 0.04466950855209243
    "T23" = 0.04695618637663403
    "G24" = 96.27664699999998
    "H24" = 288.829941
    "I24" = 0.2308894388042666
    "J24" = 0.2343807914451077
    "M24" = 47.45519633333333
    "N24" = 142.365589
    "O24" = 0.2900832029413559
    "P24" = 0.3003905575931054
    "Q24" = 4568.827185700027
    "R24" = 41119.44467130025
    "S24" = 0.06697714793367382
    "T24" = 0.07040577663130923
    "G25" = 96.27664699999998
    "H25" = 288.829941
    "I25" = 0.2308894388042666
    "J25" = 0.2343807914451077
    "M25" = 16.8400505
    "N25" = 33.680101
    "O25" = 0.102939533795646
    "P25" = 0.07106481552351887
    "Q25" = 1621.303597450673
    "R25" = 9727.82158470404
    "S25" = 0.02376765118884953
    "T25" = 0.01665622770630293
    "G26" = 96.27664699999998
    "H26" = 288.829941
    "I26" = 0.2308894388042666
    "J26" = 0.2343807914451077
    "M26" = 49.92347333333333
    "N26" = 149.77042
    "O26" = 0.3051712386725145
    "P26" = 0.3160147075621876
    "Q26" = 4806.464619127246
    "R26" = 43258.18157214521
    "S26" = 0.07046081603629976
    "T26" = 0.07406777726671977
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
